$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the date when there was no acceleration for a given speed
# (dodanie daty kiedy nie ma przyspieszenia z danej predkosci)
$ws.Range("A2").Value = "27.01.2025"
$ws.Range("A3").Value = "27.01.2025"
